$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pull_subscription")

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "pull-sub-11"
$ws.Cells.Item(7, 3).Value = "pull-sub-11"
$ws.Cells.Item(7, 4).Value = "topic-8"
$ws.Cells.Item(7, 5).Value = 120
$ws.Cells.Item(7, 6).Value = "2700s"
$ws.Cells.Item(7, 7).Value = 6
